$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-23: update date value from 45172 (2023-09-03)
# to 45175 (2023-09-06)
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
